# Update "想去人数" (wanted-to-go count) figures in column F for the
# "展览" and "全部类型" worksheets (these two sheets mirror each other).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 2820
    4  = 94
    5  = 6678
    6  = 1560
    9  = 39
    10 = 95
    12 = 131
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
